$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.293.06"
$ws.Range("E2").Value = "  +0.32%  "

# Row 3
$ws.Range("D3").Value = "1.589.35"
$ws.Range("E3").Value = "  +0.43%  "

# Row 4
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$ws.Range("D5").Value = "'212.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.16%  "

# Row 6
$ws.Range("E6").Value = "  +1.07%  "

# Row 7
$ws.Range("E7").Value = "  -0.16%  "

# Row 8
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("E9").Value = "  -0.26%  "

# Row 10
$ws.Range("D10").Value = "'19.35"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.78%  "

# Row 11
$ws.Range("E11").Value = "  +0.08%  "

# Row 12
$ws.Range("D12").Value = "1.813.12"
$ws.Range("E12").Value = "  +0.46%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.634.70"
$ws.Range("E13").Value = "  +2.41%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'4.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.24%  "

# Row 15
$ws.Range("D15").Value = "'0.520"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.89%  "

# Row 16
$ws.Range("E16").Value = "  -0.12%  "

# Row 17
$ws.Range("D17").Value = "26.309.61"
$ws.Range("E17").Value = "  +0.41%  "

# Row 18
$ws.Range("E18").Value = "  -0.47%  "

# Row 19
$ws.Range("D19").Value = "'7.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.79%  "

# Row 20
$ws.Range("D20").Value = "'213.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.13%  "

# Row 21
$ws.Range("E21").Value = "  -0.19%  "

# Row 22
$ws.Range("D22").Value = "'4.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.73%  "

# Row 23
$ws.Range("D23").Value = "'8.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.23%  "

# Row 24
$ws.Range("E24").Value = "  -2.73%  "

# Row 25
$ws.Range("D25").Value = "'145.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.12%  "

# Row 26
$ws.Range("E26").Value = "  -0.12%  "

# Row 27
$ws.Range("E27").Value = "  +0.10%  "

# Row 28
$ws.Range("E28").Value = "  -0.52%  "

# Row 29
$ws.Range("E29").Value = "  -0.27%  "

# Row 30
$ws.Range("E30").Value = "  -0.72%  "

# Row 31
$ws.Range("D31").Value = "'1.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.99%  "

# Row 32
$ws.Range("E32").Value = "  -0.32%  "

# Row 33
$ws.Range("E33").Value = "  +0.79%  "

# Row 34
$ws.Range("D34").Value = "1.341.12"
$ws.Range("E34").Value = "  +4.55%  "

# Row 35
$ws.Range("E35").Value = "  -0.91%  "

# Row 36
$ws.Range("D36").Value = "'0.599"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.58%  "

# Row 37
$ws.Range("E37").Value = "  -0.02%  "

# Row 38
$ws.Range("E38").Value = "  -0.03%  "

# Row 39
$ws.Range("E39").Value = "  -13.43%  "

# Row 40
$ws.Range("D40").Value = "'0.816"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.24%  "

# Row 41
$ws.Range("D41").Value = "'5.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.60%  "

# Row 42
$ws.Range("E42").Value = "  -0.17%  "

# Row 43
$ws.Range("E43").Value = "  +0.30%  "

# Row 44
$ws.Range("E44").Value = "  -0.57%  "

# Row 45
$ws.Range("D45").Value = "1.725.00"

# Row 46
$ws.Range("D46").Value = "'61.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.03%  "

# Row 47
$ws.Range("D47").Value = "'88.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.96%  "

# Row 48
$ws.Range("E48").Value = "  -4.17%  "

# Row 49
$ws.Range("E49").Value = "  -2.66%  "

# Row 50
$ws.Range("E50").Value = "  -0.71%  "

# Row 51
$ws.Range("E51").Value = "  -0.45%  "
